$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row numbers (1-based Excel rows) in column E whose "estado" value must be
# toggled between "Nada" and "Trabaja", per the target diff.
$rows = @(3,4,5,21,22,23,26,28,31,39,43,45,46,47,50,64,65,66,67,70,73,76,77,80,81,85,86,90,91,93,96,102,110,113,114,121,122,123,126,131,132,133,134,135,136,138,140,141,156,158,161,167,168,169,172,173,174,183,184,185,187,188,203,204,207,209,210,213,218,219,220,224,227,231,234,235,250,252,253,255,258,261,266,267,269,271,273,275,280,294,296,297,299,301,302,303,306,310,311,313,315,320,324,325,332,340,342,343,344,346,347,348,349,350,352,356,358,361,362,364,366,368)

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 5)   # Column E = 5
    $current = $cell.Value2
    if ($current -eq "Nada") {
        $cell.Value = "Trabaja"
    } elseif ($current -eq "Trabaja") {
        $cell.Value = "Nada"
    }
}
